$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the melatonin table value for E15 (No. Hours) from 11 to 19
# (latest additional citations row - totals below recalc automatically)
$ws.Range("E15").Value = 19

# Reset the saved view state: scroll back to the top-left of the sheet
# (drops the stale topLeftCell="A13") and clear the lingering F15
# selection left over from the previous edit session.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()
